# Applies the PNAD 2009 "seguranca" correction:
#  - removes the two blank "header-like" label rows that had no data
#    ("situação do domicílio" and "grandes regiões e unidades da federação")
#  - fixes the 2nd header row so that the previously "unnamed" columns now
#    read "total" (reusing the same text already present in column C)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 ("grandes regiões e unidades da federação") must be removed first so
# that row 5 still refers to "situação do domicílio" when it is deleted.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()

# Fix header row 2: columns B and F used to hold "unnamed: 1_level_1" and
# "unnamed: 5_level_1" placeholders; they now simply read "total".
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
